$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Authors (column E) with refreshed data ---
$ws.Range("E3").Value = '[Fei%Zhou%NULL%0,   Ting%Yu%NULL%0,   Ronghui%Du%NULL%0,   Guohui%Fan%NULL%0,   Ying%Liu%NULL%0,   Zhibo%Liu%NULL%0,   Jie%Xiang%NULL%0,   Yeming%Wang%NULL%0,   Bin%Song%NULL%0,   Xiaoying%Gu%NULL%0,   Lulu%Guan%NULL%0,   Yuan%Wei%NULL%0,   Hui%Li%NULL%0,   Xudong%Wu%NULL%0,   Jiuyang%Xu%NULL%0,   Shengjin%Tu%NULL%0,   Yi%Zhang%NULL%0,   Hua%Chen%NULL%0,   Bin%Cao%NULL%0]'
$ws.Range("E4").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,  Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%3,  Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%3,  Subert-Salas%Lizandra%coreGivesNoEmail%3]'
$ws.Range("E5").Value = '[Ning%Tang%NULL%0,   Dengju%Li%NULL%0,   Xiong%Wang%NULL%0,   Ziyong%Sun%NULL%0]'
$ws.Range("E6").Value = '[Rong%Qu%NULL%1,   Yun%Ling%NULL%0,   Yi‐hui‐zhi%Zhang%NULL%1,   Li‐ya%Wei%NULL%1,   Xiao%Chen%NULL%1,   Xu‐mian%Li%NULL%1,   Xuan‐yong%Liu%NULL%1,   Han‐mian%Liu%NULL%1,   Zhi%Guo%guozhi77@126.com%1,   Hua%Ren%renhua2009@hotmail.com%2,   Hua%Ren%renhua2009@hotmail.com%0,   Qiang%Wang%wangqiang@wust.edu.cn%0]'
$ws.Range("E7").Value = '[Qiurong%Ruan%NULL%0,   Kun%Yang%NULL%4,   Kun%Yang%NULL%0,   Wenxia%Wang%NULL%4,   Wenxia%Wang%NULL%0,   Lingyu%Jiang%NULL%4,   Lingyu%Jiang%NULL%0,   Jianxin%Song%songsingsjx@sina.com%4,   Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E8").Value = '[Chaolin%Huang%NULL%0,   Yeming%Wang%NULL%0,   Xingwang%Li%NULL%0,   Lili%Ren%NULL%0,   Jianping%Zhao%NULL%0,   Yi%Hu%NULL%0,   Li%Zhang%NULL%0,   Guohui%Fan%NULL%0,   Jiuyang%Xu%NULL%0,   Xiaoying%Gu%NULL%0,   Zhenshun%Cheng%NULL%0,   Ting%Yu%NULL%0,   Jiaan%Xia%NULL%0,   Yuan%Wei%NULL%0,   Wenjuan%Wu%NULL%0,   Xuelei%Xie%NULL%0,   Wen%Yin%NULL%0,   Hui%Li%NULL%0,   Min%Liu%NULL%0,   Yan%Xiao%NULL%0,   Hong%Gao%NULL%0,   Li%Guo%NULL%0,   Jungang%Xie%NULL%0,   Guangfa%Wang%NULL%0,   Rongmeng%Jiang%NULL%0,   Zhancheng%Gao%NULL%0,   Qi%Jin%NULL%0,   Jianwei%Wang%wangjw28@163.com%0,   Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("E9").Value = '[Yingxia%Liu%yingxialiu@hotmail.com%0,   Yang%Yang%NULL%0,   Cong%Zhang%NULL%0,   Fengming%Huang%NULL%0,   Fuxiang%Wang%NULL%0,   Jing%Yuan%NULL%0,   Zhaoqin%Wang%NULL%0,   Jinxiu%Li%NULL%0,   Jianming%Li%NULL%0,   Cheng%Feng%NULL%0,   Zheng%Zhang%NULL%0,   Lifei%Wang%NULL%0,   Ling%Peng%NULL%0,   Li%Chen%NULL%0,   Yuhao%Qin%NULL%0,   Dandan%Zhao%NULL%0,   Shuguang%Tan%NULL%0,   Lu%Yin%NULL%0,   Jun%Xu%NULL%0,   Congzhao%Zhou%NULL%0,   Chengyu%Jiang%jiang@pumc.edu.cn%0,   Lei%Liu%liulei3322@aliyun.com%0]'
$ws.Range("E11").Value = '[Yong%Gao%NULL%0,   Tuantuan%Li%NULL%1,   Mingfeng%Han%NULL%0,   Xiuyong%Li%NULL%1,   Dong%Wu%NULL%1,   Yuanhong%Xu%NULL%1,   Yulin%Zhu%NULL%1,   Yan%Liu%NULL%1,   Xiaowu%Wang%wangxiaowu19880218@126.com%1,   Linding%Wang%wanglinding@ahmu.edu.cn%1]'
$ws.Range("E12").Value = '[Wei-jie%Guan%NULL%0,   Zheng-yi%Ni%NULL%0,   Zheng-yi%Ni%NULL%0,   Yu%Hu%NULL%0,   Wen-hua%Liang%NULL%0,   Chun-quan%Ou%NULL%0,   Jian-xing%He%NULL%0,   Lei%Liu%NULL%0,   Hong%Shan%NULL%0,   Chun-liang%Lei%NULL%0,   David S.C.%Hui%NULL%0,   Bin%Du%NULL%0,   Lan-juan%Li%NULL%0,   Guang%Zeng%NULL%0,   Kwok-Yung%Yuen%NULL%0,   Ru-chong%Chen%NULL%0,   Chun-li%Tang%NULL%0,   Tao%Wang%NULL%0,   Ping-yan%Chen%NULL%0,   Jie%Xiang%NULL%0,   Shi-yue%Li%NULL%0,   Jin-lin%Wang%NULL%0,   Zi-jing%Liang%NULL%0,   Yi-xiang%Peng%NULL%0,   Li%Wei%NULL%0,   Yong%Liu%NULL%0,   Ya-hua%Hu%NULL%0,   Peng%Peng%NULL%0,   Jian-ming%Wang%NULL%0,   Ji-yang%Liu%NULL%0,   Zhong%Chen%NULL%0,   Gang%Li%NULL%0,   Zhi-jian%Zheng%NULL%0,   Shao-qin%Qiu%NULL%0,   Jie%Luo%NULL%0,   Chang-jiang%Ye%NULL%0,   Shao-yong%Zhu%NULL%0,   Nan-shan%Zhong%NULL%0]'
$ws.Range("E14").Value = '[Kunhua%Li%NULL%0,   Jiong%Wu%NULL%1,   Faqi%Wu%NULL%1,   Dajing%Guo%NULL%1,   Linli%Chen%NULL%1,   Zheng%Fang%NULL%1,   Chuanming%Li%NULL%1]'
$ws.Range("E15").Value = '[Wei%Liu%NULL%0,   Zhao-Wu%Tao%NULL%3,   Lei%Wang%NULL%3,   Ming-Li%Yuan%NULL%0,   Kui%Liu%NULL%0,   Ling%Zhou%NULL%0,   Shuang%Wei%NULL%0,   Yan%Deng%NULL%0,   Jing%Liu%NULL%5,   Hui-Guo%Liu%NULL%0,   Ming%Yang%NULL%4,   Yi%Hu%NULL%0,   Pei-Fang%Wei%NULL%0,   Pei-Fang%Wei%NULL%0]'
$ws.Range("E16").Value = '[Chuan%Qin%NULL%2,   Luoqi%Zhou%NULL%1,   Ziwei%Hu%NULL%1,   Shuoqi%Zhang%NULL%1,   Sheng%Yang%NULL%1,   Yu%Tao%NULL%0,   Cuihong%Xie%NULL%0,   Ke%Ma%NULL%0,   Ke%Shang%NULL%1,   Wei%Wang%NULL%0,   Dai-Shi%Tian%tiands@tjh.tjmu.edu.cn%1]'
$ws.Range("E17").Value = '[Dan%Sun%NULL%0,   Hui%Li%NULL%0,   Xiao-Xia%Lu%NULL%1,   Han%Xiao%NULL%0,   Jie%Ren%NULL%1,   Fu-Rong%Zhang%792523496@qq.com%1,   Zhi-Sheng%Liu%liuzsc@126.com%1]'
$ws.Range("E18").Value = '[Suxin%Wan%NULL%0,   Yi%Xiang%NULL%0,   Wei%Fang%NULL%0,   Yu%Zheng%yuzheng1@cdutcm.edu.cn%0,   Boqun%Li%279685211@qq.com%0,   Boqun%Li%279685211@qq.com%0,   Yanjun%Hu%huyanjun@163.com%0,   Chunhui%Lang%NULL%0,   Daoqiu%Huang%NULL%0,   Qiuyan%Sun%NULL%0,   Yan%Xiong%NULL%0,   Xia%Huang%NULL%0,   Jinglong%Lv%NULL%0,   Yaling%Luo%NULL%0,   Li%Shen%NULL%0,   Haoran%Yang%NULL%0,   Gu%Huang%NULL%0,   Ruishan%Yang%NULL%0]'
$ws.Range("E19").Value = '[Zhongliang%Wang%NULL%0,   Bohan%Yang%NULL%0,   Qianwen%Li%NULL%0,   Lu%Wen%NULL%0,   Ruiguang%Zhang%zrg27@163.com%0]'
$ws.Range("E21").Value = '[Xiaobo%Yang%NULL%0,   Yuan%Yu%NULL%0,   Jiqian%Xu%NULL%0,   Huaqing%Shu%NULL%0,   Jia''an%Xia%NULL%0,   Hong%Liu%NULL%0,   Yongran%Wu%NULL%0,   Lu%Zhang%NULL%0,   Zhui%Yu%NULL%0,   Minghao%Fang%NULL%0,   Ting%Yu%NULL%0,   Yaxin%Wang%NULL%0,   Shangwen%Pan%NULL%0,   Xiaojing%Zou%NULL%0,   Shiying%Yuan%NULL%0,   You%Shang%NULL%0]'
$ws.Range("E22").Value = '[Akdis%Cezmi A%coreGivesNoEmail%5,  Cao%Yi-yuan%coreGivesNoEmail%4,  Dong%Xiang%coreGivesNoEmail%4,  Gao%Ya-dong%coreGivesNoEmail%4,  Yan%You-qin%coreGivesNoEmail%4,  Yang%Yi-bin%coreGivesNoEmail%4,  Yuan%Ya-dong%coreGivesNoEmail%4,  Zhang%Jin-jin%coreGivesNoEmail%4]'

# --- Update ID (F) / ID Format (G) columns for rows 4 and 22 ---
$ws.Range("F4").Value = 'not found'
$ws.Range("G4").Value = 'N/A'
$ws.Range("F22").Value = 'not found'
$ws.Range("G22").Value = 'N/A'

# --- Add new column I: "Other found locations" ---
$ws.Range("I1").Value = 'Other found locations'
$ws.Range("I2").Value = ''
$ws.Range("I3").Value = '_PMC_elsevier'
$ws.Range("I4").Value = ''
$ws.Range("I5").Value = '_PMC_elsevier'
$ws.Range("I6").Value = '_PMC'
$ws.Range("I7").Value = '_PMC_Springer'
$ws.Range("I8").Value = '_PMC_elsevier'
$ws.Range("I9").Value = '_PMC_Springer'
$ws.Range("I10").Value = '_PMC'
$ws.Range("I11").Value = '_PMC'
$ws.Range("I12").Value = '_PMC'
$ws.Range("I13").Value = ''
$ws.Range("I14").Value = '_PMC'
$ws.Range("I15").Value = '_PMC'
$ws.Range("I16").Value = '_PMC'
$ws.Range("I17").Value = '_PMC_Springer'
$ws.Range("I18").Value = '_PMC'
$ws.Range("I19").Value = '_PMC'
$ws.Range("I20").Value = ''
$ws.Range("I21").Value = '_PMC_elsevier'
$ws.Range("I22").Value = ''
